$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# WorldID and ServerID columns (B and C) change their type from "int" to "int16_t"
$ws.Range("B3").Value = "int16_t"
$ws.Range("C3").Value = "int16_t"

# Reflect the last active cell as C3 (matches where the edit took place)
$ws.Range("C3").Select()
